# Fixed errors of dataprovider
# - Insert a new leading column "Automation Test Script ID" (ATC00x) before the
#   existing TC_ID column.
# - Correct the typo "2 Travellers(s)" -> "2 Traveller(s)" in the TC_003 row.
# - Update the active selection to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the front; this shifts all existing columns (and their
# formatting/column widths) one position to the right.
$ws.Columns("A:A").Insert()

# New header cell, matching the bold style used by the other header cells.
$ws.Range("A1").Value = "Automation Test Script ID"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Automation Test Script IDs for each test case row.
$ws.Range("A2").Value = "ATC001"
$ws.Range("A3").Value = "ATC002"
$ws.Range("A4").Value = "ATC003"
$ws.Range("A5").Value = "ATC004"
$ws.Range("A6").Value = "ATC005"
$ws.Range("A7").Value = "ATC006"
$ws.Range("A8").Value = "ATC007"
$ws.Range("A9").Value = "ATC008"
$ws.Range("A10").Value = "ATC009"
$ws.Range("A11").Value = "ATC010"

# Fix the "2 Travellers(s)" typo (row for TC_003, now column G after the insert).
$ws.Range("G4").Value = "2 Traveller(s)"

# Update the selected cell to reflect the author's final cursor position.
$ws.Range("C5").Select()
